$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift each date in F2:F7 forward by one day (add 1 to underlying serial value)
for ($row = 2; $row -le 7; $row++) {
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value = $cell.Value2 + 1
}
